# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The underlying worker/employee rows (B16:J22) are re-sorted into a new
# order. Every row keeps its own formatting (s="..." style ids, borders,
# etc.) - only the data (document number, name, period, valor mora,
# salario basico) moves between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (document type, document number, name, period, valor mora, salario basico)
$data = @(
    @("CC", "1128046927", "EDER JULIAN ARCHBOLD SALCEDO",     "1809", 1200,  900000),
    @("CC", "1143347543", "ALEMIS VILLARREAL ANGULO",         "1809", 1053,  790000),
    @("CC", "1047421288", "JOHNNY FERNANDO REINA BOLIVAR",    "1903", 4417,  828116),
    @("CC", "73205202",   "HAROLD JEYSON HERRERA SAYAVEDRA",  "1905", 1893,  1420000),
    @("CC", "1143401657", "VERONICA RAMOS LARA",              "1905", 1104,  877803),
    @("CC", "73214033",   "ENOTH ENRIQUE GARCIA YEPEZ",       "2006", 35112, 877803),
    @("CC", "1143363534", "CARLOS DE JESUS MARTINEZ VILORIA", "2412", 24800, 1550000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 2).Value = $rowData[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rowData[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rowData[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rowData[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rowData[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $rowData[5]   # G - Salario Basico
}
